$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = '57.067.69'
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = '  +7.13%  '
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = '3.238.99'
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = '  +2.93%  '
$ws.Cells.Item(4,5).Value = '  +0.01%  '
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '394.29'
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = '  -1.02%  '
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '106.97'
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = '  +0.54%  '
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '3.235.71'
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = '  +2.77%  '
$ws.Cells.Item(8,5).Value = '  +3.44%  '
$ws.Cells.Item(9,5).Value = '  -0.04%  '
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '0.616'
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = '  +1.40%  '
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '38.85'
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = '  -0.03%  '
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '0.0977'
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = '  +12.16%  '
$ws.Cells.Item(13,5).Value = '  +1.68%  '
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '3.753.39'
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = '  +2.95%  '
$ws.Cells.Item(15,5).Value = '  +2.14%  '
$ws.Cells.Item(16,5).Value = '  -0.18%  '
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '3.240.05'
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value = '  +2.74%  '
$ws.Cells.Item(18,5).Value = '  -1.70%  '
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '10.99'
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = '  +1.70%  '
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '56.838.62'
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = '  +6.79%  '
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '3.33'
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = '  +1.44%  '
$ws.Cells.Item(22,5).Value = '  +7.98%  '
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '12.95'
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = '  +0.42%  '
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '298.10'
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = '  +9.91%  '
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '73.59'
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = '  +3.42%  '
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '3.13'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = '  -2.72%  '
$ws.Cells.Item(27,2).Value = 'LEO'
$ws.Cells.Item(27,3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '4.38'
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = '  +3.30%  '
$ws.Cells.Item(28,2).Value = 'EthereumClassic'
$ws.Cells.Item(28,3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '27.83'
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value = '  +0.82%  '
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '7.70'
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value = '  -4.25%  '
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '7.23'
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value = '  -3.54%  '
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '0.999'
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value = '  -0.01%  '
$ws.Cells.Item(33,5).Value = '  -0.80%  '
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '10.93'
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = '  -0.80%  '
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '37.21'
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value = '  -0.05%  '
$ws.Cells.Item(36,5).Value = '  -2.20%  '
$ws.Cells.Item(37,5).Value = '  +1.13%  '
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '51.60'
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value = '  +2.28%  '
$ws.Cells.Item(39,2).Value = 'Stacks'
$ws.Cells.Item(39,3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '3.10'
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = '  +12.88%  '
$ws.Cells.Item(40,2).Value = 'LidoDAOToken'
$ws.Cells.Item(40,3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '3.51'
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = '  -1.15%  '
$ws.Cells.Item(41,2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(41,3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.999'
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = '  -0.09%  '
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '133.82'
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = '  +2.90%  '
$ws.Cells.Item(43,5).Value = '  +0.16%  '
$ws.Cells.Item(44,5).Value = '  +1.92%  '
$ws.Cells.Item(45,5).Value = '  -4.98%  '
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '16.89'
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = '  -2.71%  '
$ws.Cells.Item(47,5).Value = '  -3.90%  '
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '21.90'
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = '  -1.84%  '
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '2.144.33'
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = '  +2.47%  '
$ws.Cells.Item(50,5).Value = '  -0.50%  '
$ws.Cells.Item(51,5).Value = '  +24.89%  '
